$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('B1').Value = 'Sample description'
$ws.Range('E1').Value = 'Iron II'
$ws.Range('F1').Value = 'MN II'
$ws.Range('G1').Value = 'cyanid'
$ws.Range('H1').Value = 'phosphor'
$ws.Range('I1').Value = 'benzene'
$ws.Range('J1').Value = 'toluene'
$ws.Range('K1').Value = 'ethylbenzene'
$ws.Range('L1').Value = 'o-xylene'
$ws.Range('M1').Value = '(m+p)-xylene'
$ws.Range('N1').Value = 'sum xylenes (factor 0.7)'
$ws.Range('O1').Value = 'total BTEX (factor 0.7)'
$ws.Range('P1').Value = 'naphthalene'
$ws.Range('Q1').Value = 'phenol'
$ws.Range('U1').Value = 'som cresols'
$ws.Range('V1').Value = '2-ethylphenol'
$ws.Range('W1').Value = '3-ethylphenol'
$ws.Range('X1').Value = '2,4-dimethylphenol'
$ws.Range('Y1').Value = '2,5-dimethylphenol'
$ws.Range('Z1').Value = '3,5+2,3-dimethylphenol+4-ethylphenol'
$ws.Range('AA1').Value = '2,6-dimethylphenol'
$ws.Range('AB1').Value = '3,4-dimethylphenol'
$ws.Range('AC1').Value = 'som C2-alkylphenolen'
$ws.Range('AD1').Value = '2,3,5-trimethylphenol'
$ws.Range('AE1').Value = '3,4,5-trimethylphenol'
$ws.Range('AF1').Value = '2-isopropylphenol'
$ws.Range('AG1').Value = 'som C3-alkylphenolen'
$ws.Range('AI1').Value = 'p-(tert)butylphenol'
$ws.Range('AJ1').Value = 'som C4-alkylphenolen'
$ws.Range('AK1').Value = 'naphthalene'
$ws.Range('AL1').Value = 'acenaphthylene'
$ws.Range('AM1').Value = 'acenaphtene'
$ws.Range('AN1').Value = 'fluorene'
$ws.Range('AO1').Value = 'phenanthrene'
$ws.Range('AP1').Value = 'anthracene'
$ws.Range('AQ1').Value = 'fluoranthene'
$ws.Range('AR1').Value = 'pyrene'
$ws.Range('AT1').Value = 'chrysene'
$ws.Range('AU1').Value = 'benzo(b)fluoranthene'
$ws.Range('AV1').Value = 'benzo(k)fluoranthene'
$ws.Range('AW1').Value = 'benzo(a)pyrene'
$ws.Range('AX1').Value = 'dibenz(a,h)anthracene'
$ws.Range('AY1').Value = 'benzo(g,h,i)perylene'
$ws.Range('AZ1').Value = 'indeno(1,2,3-cd)pyrene'
$ws.Range('BA1').Value = 'sum PAH (16 EPA)'
$ws.Range('BB1').Value = 'sum PAH (VROM) (factor 0.7)'
$ws.Range('BC1').Value = 'fraction C10-C12'
$ws.Range('BD1').Value = 'fraction C12-C22'
$ws.Range('BE1').Value = 'fraction C22-C30'
$ws.Range('BF1').Value = 'fraction C30-C40'
$ws.Range('BG1').Value = 'total oil C10 - C40'
$ws.Range('BI1').Value = 'nitrite'
$ws.Range('BJ1').Value = 'nitrite - N'
$ws.Range('BK1').Value = 'nitrate'
$ws.Range('BL1').Value = 'nitrate - N'
$ws.Range('BM1').Value = 'sulphates'
$ws.Range('BN1').Value = 'Oxygen'
